# Applies weekly fruit/vegetable price update: reshuffles per-row data
# (date, volume, prices, unit, origin, price/kg, kg-or-units) across rows 2-15
# per the source commit, while keeping market/category columns intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 10
$ws.Range("D2").Value = 44284
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("N2").Value = '$/caja 60 unidades'
$ws.Range("O2").Value = 'Provincia de Limarí'
$ws.Range("P2").Value = 167
$ws.Range("Q2").Value = 60

# Row 3 <- source row 12
$ws.Range("D3").Value = 44243
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 10375
$ws.Range("N3").Value = '$/caja 60 unidades'
$ws.Range("O3").Value = 'Provincia de Quillota'
$ws.Range("P3").Value = 173
$ws.Range("Q3").Value = 60

# Row 4 <- source row 13
$ws.Range("D4").Value = 44333
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 10400
$ws.Range("N4").Value = '$/caja 60 unidades'
$ws.Range("O4").Value = 'Provincia de Limarí'
$ws.Range("P4").Value = 173
$ws.Range("Q4").Value = 60

# Row 5 <- source row 7
$ws.Range("D5").Value = 44291
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("N5").Value = '$/caja 60 unidades'
$ws.Range("O5").Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 150
$ws.Range("Q5").Value = 60

# Row 6 <- source row 5
$ws.Range("D6").Value = 44186
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("N6").Value = '$/caja 60 unidades'
$ws.Range("O6").Value = 'Provincia de Limarí'
$ws.Range("P6").Value = 117
$ws.Range("Q6").Value = 60

# Row 7 <- source row 15
$ws.Range("D7").Value = 44200
$ws.Range("J7").Value = 10
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 9000
$ws.Range("N7").Value = '$/caja 60 unidades'
$ws.Range("O7").Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 150
$ws.Range("Q7").Value = 60

# Row 8 <- source row 14
$ws.Range("D8").Value = 45030
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 6000
$ws.Range("N8").Value = '$/caja 50 unidades'
$ws.Range("O8").Value = 'Región de Arica y Parinacota'
$ws.Range("P8").Value = 120
$ws.Range("Q8").Value = 50

# Row 9 <- source row 4
$ws.Range("D9").Value = 45001
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("N9").Value = '$/caja 60 unidades'
$ws.Range("O9").Value = 'Provincia de Limarí'
$ws.Range("P9").Value = 167
$ws.Range("Q9").Value = 60

# Row 10 <- source row 11
$ws.Range("D10").Value = 44315
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("N10").Value = '$/caja 60 unidades'
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("P10").Value = 167
$ws.Range("Q10").Value = 60

# Row 11 <- source row 8
$ws.Range("D11").Value = 44405
$ws.Range("J11").Value = 45
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 9000
$ws.Range("M11").Value = 9000
$ws.Range("N11").Value = '$/caja 50 unidades'
$ws.Range("O11").Value = 'Provincia de Quillota'
$ws.Range("P11").Value = 180
$ws.Range("Q11").Value = 50

# Row 12 <- source row 2
$ws.Range("D12").Value = 44585
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 11000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 11000
$ws.Range("N12").Value = '$/caja 60 unidades'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 183
$ws.Range("Q12").Value = 60

# Row 13 <- source row 6
$ws.Range("D13").Value = 44312
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("N13").Value = '$/caja 60 unidades'
$ws.Range("O13").Value = 'Provincia de Limarí'
$ws.Range("P13").Value = 167
$ws.Range("Q13").Value = 60

# Row 14 <- source row 9
$ws.Range("D14").Value = 44277
$ws.Range("J14").Value = 25
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("N14").Value = '$/caja 60 unidades'
$ws.Range("O14").Value = 'Provincia de Limarí'
$ws.Range("P14").Value = 167
$ws.Range("Q14").Value = 60

# Row 15 <- source row 3
$ws.Range("D15").Value = 44179
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("N15").Value = '$/caja 60 unidades'
$ws.Range("O15").Value = 'Provincia de Limarí'
$ws.Range("P15").Value = 117
$ws.Range("Q15").Value = 60

